# Refresh the cryptos price/volume table (GitHub Actions daily update).
# D-column price cells are forced to Text ("@") before assignment and the
# style is reset to "Normal" afterwards so values like "240.10" / "159.00"
# are kept verbatim as strings instead of being coerced into numbers
# (which would silently drop the trailing zero, e.g. 240.10 -> 240.1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.168.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.246.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.18%  "
$ws.Range("E11").Value = "  +6.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0793"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("E13").Value = "  +3.16%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.595.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.219.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.105.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.48%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.17%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.10%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.61%  "
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.54%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.05%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("E40").Value = "  +5.26%  "
$ws.Range("E41").Value = "  +6.20%  "
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.074.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.84%  "
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.467.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.58%  "
